# Auto-generated edit script: updates crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "42.794.17"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.298.86"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "305.95"
$ws.Range("E5").Value = "  +1.86%  "
Set-TextValue $ws.Range("D6") "96.74"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("E8").Value = "  +0.06%  "
Set-TextValue $ws.Range("D10") "35.53"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("E11").Value = "  +0.22%  "
Set-TextValue $ws.Range("D12") "18.31"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "2.657.86"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "2.303.63"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "42.749.29"
$ws.Range("E18").Value = "  -0.45%  "
Set-TextValue $ws.Range("D19") "13.03"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  -1.19%  "
Set-TextValue $ws.Range("D21") "6.02"
$ws.Range("E21").Value = "  -1.72%  "
Set-TextValue $ws.Range("D22") "67.33"
$ws.Range("E22").Value = "  -1.28%  "
Set-TextValue $ws.Range("D23") "236.10"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  -2.32%  "
Set-TextValue $ws.Range("D25") "2.46"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.00%  "
Set-TextValue $ws.Range("D27") "4.03"
$ws.Range("E27").Value = "  +0.32%  "
Set-TextValue $ws.Range("D28") "25.37"
$ws.Range("E28").Value = "  +1.16%  "
Set-TextValue $ws.Range("D29") "166.70"
$ws.Range("E29").Value = "  +2.57%  "
Set-TextValue $ws.Range("D30") "2.07"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("E31").Value = "  -0.92%  "
Set-TextValue $ws.Range("D32") "33.28"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  -0.63%  "
Set-TextValue $ws.Range("D38") "0.0691"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("E39").Value = "  -1.00%  "
Set-TextValue $ws.Range("D40") "1.75"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").Value = "2.002.74"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").Value = "  -2.23%  "
Set-TextValue $ws.Range("D45") "18.29"
$ws.Range("E45").Value = "  +4.37%  "
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("E47").Value = "  -5.29%  "
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D49") "53.72"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.526.43"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D51") "2.82"
$ws.Range("E51").Value = "  +3.48%  "
